$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Save" column, reuse the style from G1 (bold/bordered header style)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (rows 2-13)
$saveValues = @(0, 1, 0, 0, 0, 1, 0, 0, 0, 1, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
